# Update species names in the LHS settings input workbook:
#   BV1 -> NO, LB1 -> Li, P -> oLB
# Rewrites the "name" column (column A) on both the "parameters" and
# "initial_conditions" sheets so every k_{grow}-*, \alpha_{*->*} and bare
# species label reflects the new naming scheme, then restores the
# view/selection state (active sheet, zoom, selected cell) recorded for
# each sheet.

$wb = $excel.ActiveWorkbook

$params = $wb.Worksheets.Item("parameters")
$ic = $wb.Worksheets.Item("initial_conditions")

# --- "parameters" sheet: column A (rows 2-13) ---
$params.Range("A2").Value = "k_{grow}-NO"
$params.Range("A3").Value = "k_{grow}-Li"
$params.Range("A4").Value = "k_{grow}-oLB"

$params.Range("A5").Value  = "\alpha_{NO->NO}"
$params.Range("A6").Value  = "\alpha_{NO->Li}"
$params.Range("A7").Value  = "\alpha_{NO->oLB}"
$params.Range("A8").Value  = "\alpha_{Li->NO}"
$params.Range("A9").Value  = "\alpha_{Li->Li}"
$params.Range("A10").Value = "\alpha_{Li->oLB}"
$params.Range("A11").Value = "\alpha_{oLB->NO}"
$params.Range("A12").Value = "\alpha_{oLB->Li}"
$params.Range("A13").Value = "\alpha_{oLB->oLB}"

# --- "initial_conditions" sheet: column A (rows 2-4) ---
$ic.Range("A2").Value = "NO"
$ic.Range("A3").Value = "Li"
$ic.Range("A4").Value = "oLB"

# --- view state: active sheet / zoom / selection ---
$ic.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 174
$ic.Range("A4").Select() | Out-Null

$params.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 208
$params.Range("A14").Select() | Out-Null
